$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$r = $ws.Range("D2")
$r.NumberFormat = "@"
$r.Value = "64.612.06"
$r.ClearFormats()
$ws.Range("E2").Value = "  -0.43%  "
$r = $ws.Range("D3")
$r.NumberFormat = "@"
$r.Value = "3.423.03"
$r.ClearFormats()
$ws.Range("E3").Value = "  -1.08%  "
$ws.Range("E4").Value = "  +0.00%  "
$r = $ws.Range("D5")
$r.NumberFormat = "@"
$r.Value = "573.19"
$r.ClearFormats()
$ws.Range("E5").Value = "  -0.50%  "
$r = $ws.Range("D6")
$r.NumberFormat = "@"
$r.Value = "156.92"
$r.ClearFormats()
$ws.Range("E6").Value = "  -2.19%  "
$r = $ws.Range("D7")
$r.NumberFormat = "@"
$r.Value = "0.622"
$r.ClearFormats()
$ws.Range("E7").Value = "  +7.04%  "
$ws.Range("E8").Value = "  +0.07%  "
$r = $ws.Range("D9")
$r.NumberFormat = "@"
$r.Value = "3.427.73"
$r.ClearFormats()
$ws.Range("E9").Value = "  -1.01%  "
$ws.Range("E10").Value = "  -2.69%  "
$r = $ws.Range("D11")
$r.NumberFormat = "@"
$r.Value = "0.123"
$r.ClearFormats()
$ws.Range("E11").Value = "  -1.92%  "
$r = $ws.Range("D12")
$r.NumberFormat = "@"
$r.Value = "0.443"
$r.ClearFormats()
$ws.Range("E12").Value = "  +0.53%  "
$r = $ws.Range("D13")
$r.NumberFormat = "@"
$r.Value = "4.016.57"
$r.ClearFormats()
$ws.Range("E13").Value = "  -1.02%  "
$ws.Range("E14").Value = "  +0.39%  "
$ws.Range("E15").Value = "  -3.36%  "
$r = $ws.Range("D16")
$r.NumberFormat = "@"
$r.Value = "27.90"
$r.ClearFormats()
$r = $ws.Range("D17")
$r.NumberFormat = "@"
$r.Value = "64.615.24"
$r.ClearFormats()
$ws.Range("E17").Value = "  -0.47%  "
$r = $ws.Range("D18")
$r.NumberFormat = "@"
$r.Value = "3.425.39"
$r.ClearFormats()
$ws.Range("E18").Value = "  -1.98%  "
$r = $ws.Range("D19")
$r.NumberFormat = "@"
$r.Value = "6.37"
$r.ClearFormats()
$ws.Range("E19").Value = "  -0.04%  "
$r = $ws.Range("D20")
$r.NumberFormat = "@"
$r.Value = "13.99"
$r.ClearFormats()
$ws.Range("E20").Value = "  -2.20%  "
$r = $ws.Range("D21")
$r.NumberFormat = "@"
$r.Value = "378.05"
$r.ClearFormats()
$ws.Range("E21").Value = "  -2.91%  "
$ws.Range("E22").Value = "  -1.91%  "
$r = $ws.Range("D23")
$r.NumberFormat = "@"
$r.Value = "0.550"
$r.ClearFormats()
$ws.Range("E23").Value = "  +1.10%  "
$ws.Range("E24").Value = "  -0.33%  "
$r = $ws.Range("D25")
$r.NumberFormat = "@"
$r.Value = "72.57"
$r.ClearFormats()
$ws.Range("E25").Value = "  -0.92%  "
$ws.Range("E26").Value = "  -4.38%  "
$ws.Range("E27").Value = "  +6.41%  "
$ws.Range("E28").Value = "  -1.19%  "
$r = $ws.Range("D29")
$r.NumberFormat = "@"
$r.Value = "1.00"
$r.ClearFormats()
$ws.Range("E29").Value = "  +0.00%  "
$r = $ws.Range("D30")
$r.NumberFormat = "@"
$r.Value = "1.50"
$r.ClearFormats()
$ws.Range("E30").Value = "  +4.41%  "
$r = $ws.Range("D31")
$r.NumberFormat = "@"
$r.Value = "6.20"
$r.ClearFormats()
$ws.Range("E31").Value = "  -0.41%  "
$ws.Range("E32").Value = "  -0.57%  "
$r = $ws.Range("D33")
$r.NumberFormat = "@"
$r.Value = "23.16"
$r.ClearFormats()
$ws.Range("E33").Value = "  -2.34%  "
$ws.Range("E34").Value = "  +1.83%  "
$ws.Range("E35").Value = "  +6.63%  "
$r = $ws.Range("D36")
$r.NumberFormat = "@"
$r.Value = "159.47"
$r.ClearFormats()
$ws.Range("E36").Value = "  -2.27%  "
$r = $ws.Range("D37")
$r.NumberFormat = "@"
$r.Value = "1.90"
$r.ClearFormats()
$ws.Range("E37").Value = "  -0.83%  "
$ws.Range("E39").Value = "  +0.18%  "
$r = $ws.Range("D40")
$r.NumberFormat = "@"
$r.Value = "26.92"
$r.ClearFormats()
$ws.Range("E40").Value = "  -1.47%  "
$r = $ws.Range("D41")
$r.NumberFormat = "@"
$r.Value = "2.880.04"
$r.ClearFormats()
$ws.Range("E41").Value = "  -4.37%  "
$ws.Range("E42").Value = "  +1.46%  "
$r = $ws.Range("D43")
$r.NumberFormat = "@"
$r.Value = "26.68"
$r.ClearFormats()
$ws.Range("E43").Value = "  +9.23%  "
$r = $ws.Range("D44")
$r.NumberFormat = "@"
$r.Value = "0.0318"
$r.ClearFormats()
$ws.Range("E44").Value = "  +0.79%  "
$r = $ws.Range("D45")
$r.NumberFormat = "@"
$r.Value = "42.93"
$r.ClearFormats()
$ws.Range("E45").Value = "  -0.06%  "
$ws.Range("E46").Value = "  -0.23%  "
$r = $ws.Range("D47")
$r.NumberFormat = "@"
$r.Value = "321.57"
$r.ClearFormats()
$ws.Range("E47").Value = "  +5.35%  "
$ws.Range("E48").Value = "  -0.08%  "
$ws.Range("B49").Value = "dogwifhat"
$ws.Range("C49").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$r = $ws.Range("D49")
$r.NumberFormat = "@"
$r.Value = "2.20"
$r.ClearFormats()
$ws.Range("E49").Value = "  +1.68%  "
$ws.Range("B50").Value = "Stellar"
$ws.Range("C50").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$r = $ws.Range("D50")
$r.NumberFormat = "@"
$r.Value = "0.109"
$r.ClearFormats()
$ws.Range("E50").Value = "  +2.09%  "
$ws.Range("B51").Value = "SuiNetwork"
$ws.Range("C51").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$r = $ws.Range("D51")
$r.NumberFormat = "@"
$r.Value = "0.863"
$r.ClearFormats()
$ws.Range("E51").Value = "  -1.68%  "
